$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.548.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06126"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.664.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06947"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "74.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5718"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.557.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006729"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.878.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.408"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.663"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.224"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.711"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "103.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.955"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07693"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.593"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04313"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("E34").Value = "  +1.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9429"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9198"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.477"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.38%  "

$ws.Range("E39").Value = "  +6.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.840"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01462"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.01%  "

$ws.Range("E43").Value = "  +7.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3705"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("E45").Value = "  +1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05255"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.118"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "29.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.571"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.52%  "

$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
